$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the AXP row entirely (row 4 in the original sheet)
$ws.Rows("4:4").Delete()

# Insert a new "Quantity" column before the old Difference column (D)
$ws.Columns("D:D").Insert()

# Header row
$ws.Range("D1").Value = "Quantity"
$ws.Range("E1").Value = "Difference"
$ws.Range("F1").Value = "Total gain (%)"

# Make sure quantity/difference/gain columns hold text (matches source formatting)
$ws.Range("D2:F25").NumberFormat = "@"

$ws.Range("C2").Value = 178.9503
$ws.Range("D2").Value = "2.00"
$ws.Range("E2").Value = "58.60"
$ws.Range("F2").Value = "48.69 %"
$ws.Range("C3").Value = 112.68
$ws.Range("D3").Value = "1.00"
$ws.Range("E3").Value = "28.08"
$ws.Range("F3").Value = "33.19 %"
$ws.Range("C4").Value = 203.55
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "55.60"
$ws.Range("F4").Value = "37.58 %"
$ws.Range("C5").Value = 44.5664
$ws.Range("D5").Value = "6.00"
$ws.Range("E5").Value = "15.59"
$ws.Range("F5").Value = "53.82 %"
$ws.Range("C6").Value = 47649.16
$ws.Range("D6").Value = "2.00"
$ws.Range("E6").Value = "2816.16"
$ws.Range("F6").Value = "6.28 %"
$ws.Range("C7").Value = 827
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "108.00"
$ws.Range("F7").Value = "15.02 %"
$ws.Range("C8").Value = 63.8
$ws.Range("D8").Value = "3.00"
$ws.Range("E8").Value = "23.47"
$ws.Range("F8").Value = "58.19 %"
$ws.Range("C9").Value = 104
$ws.Range("D9").Value = "3.00"
$ws.Range("E9").Value = "31.52"
$ws.Range("F9").Value = "43.49 %"
$ws.Range("C10").Value = 156.355
$ws.Range("D10").Value = "1.00"
$ws.Range("E10").Value = "55.61"
$ws.Range("F10").Value = "55.21 %"
$ws.Range("C11").Value = 137.2
$ws.Range("D11").Value = "2.00"
$ws.Range("E11").Value = "-17.30"
$ws.Range("F11").Value = "-11.20 %"
$ws.Range("C12").Value = 52.01333333333333
$ws.Range("D12").Value = "4.00"
$ws.Range("E12").Value = "3.86"
$ws.Range("F12").Value = "8.01 %"
$ws.Range("C13").Value = 58.94333333333334
$ws.Range("D13").Value = "5.00"
$ws.Range("E13").Value = "12.72"
$ws.Range("F13").Value = "27.51 %"
$ws.Range("C14").Value = 93.75
$ws.Range("D14").Value = "1.00"
$ws.Range("E14").Value = "39.65"
$ws.Range("F14").Value = "73.29 %"
$ws.Range("C15").Value = 340.94
$ws.Range("D15").Value = "2.00"
$ws.Range("E15").Value = "99.50"
$ws.Range("F15").Value = "41.21 %"
$ws.Range("C16").Value = 58.59
$ws.Range("D16").Value = "1.00"
$ws.Range("E16").Value = "26.76"
$ws.Range("F16").Value = "84.07 %"
$ws.Range("C17").Value = 162.8956
$ws.Range("D17").Value = "1.00"
$ws.Range("E17").Value = "27.10"
$ws.Range("F17").Value = "19.95 %"
$ws.Range("C18").Value = 94.55
$ws.Range("D18").Value = "3.00"
$ws.Range("E18").Value = "10.39"
$ws.Range("F18").Value = "12.35 %"
$ws.Range("C19").Value = 22.14
$ws.Range("D19").Value = "3.00"
$ws.Range("E19").Value = "11.23"
$ws.Range("F19").Value = "102.93 %"
$ws.Range("C20").Value = 24.78
$ws.Range("D20").Value = "2.00"
$ws.Range("E20").Value = "-6.66"
$ws.Range("F20").Value = "-21.17 %"
$ws.Range("C21").Value = 126.36
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "66.22"
$ws.Range("F21").Value = "110.11 %"
$ws.Range("C22").Value = 24.795
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "-4.88"
$ws.Range("F22").Value = "-16.43 %"
$ws.Range("C23").Value = 53.63
$ws.Range("D23").Value = "0.03"
$ws.Range("E23").Value = "-4.14"
$ws.Range("F23").Value = "-7.17 %"
$ws.Range("C24").Value = 258.4
$ws.Range("D24").Value = "2.00"
$ws.Range("E24").Value = "110.40"
$ws.Range("F24").Value = "74.59 %"
$ws.Range("C25").Value = 49.17
$ws.Range("D25").Value = "2.00"
$ws.Range("E25").Value = "26.56"
$ws.Range("F25").Value = "117.47 %"

$ws.Range("A1").Select()
